# BDOT10k ontology cleanup: fix Polish/English typos in the "cellfie import
# data" sheet, per the commit "English labels added to BDOT10k ontology".
#
# The workbook has a single worksheet; operate on it directly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text corrections -----------------------------------------------------
# Each assignment below corresponds to one <si>/<t> fix in the diff. Column E
# holds the Polish label, column F the English label; most of the edits are
# in column E, with one English wording fix in F128.

# "auostrada" -> "autostrada"
$ws.Range("E15").Value = "autostrada"

# "zabudowa przemysłowo-składowa" -> proper Unicode hyphen (U+2010)
$ws.Range("E50").Value = "zabudowa przemysłowo‐składowa"

# "zabudowa handlowo-usługowa" -> Unicode hyphen
$ws.Range("E51").Value = "zabudowa handlowo‐usługowa"

# "pozostały grunt nieużytkowy" -> "pozostały grunt nieużytkowany"
$ws.Range("E72").Value = "pozostały grunt nieużytkowany"

# "teren przemysłowo-składowy" -> Unicode hyphen
$ws.Range("E79").Value = "teren przemysłowo‐składowy"

# "budynki zakwaterowania turystycznego, pozostałe" -> drop the space after
# the comma
$ws.Range("E85").Value = "budynki zakwaterowania turystycznego,pozostałe"

# "kultur religijnego" -> "kultu religijnego"
$ws.Range("E98").Value = "budynki przeznaczone do sprawowania kultu religijnego i czynności religijnych"

# "tor sachochodowy" -> "tor samochodowy"
$ws.Range("E122").Value = "tor samochodowy"

# "telecommunications mast" -> "telecommunication mast" (English column)
$ws.Range("F128").Value = "telecommunication mast"

# "slup energetyczny" -> "słup energetyczny"
$ws.Range("E130").Value = "słup energetyczny"

# "zespół urządzeń stacji metereologicznej" -> "...meteorologicznej"
$ws.Range("E158").Value = "zespół urządzeń stacji meteorologicznej"

# "koplania" -> "kopalnia"
$ws.Range("E174").Value = "kopalnia"

# "centrum handlowo-usługowe" -> Unicode hyphen + trailing newline
$ws.Range("E185").Value = "centrum handlowo‐usługowe`n"

# "ośrodek sportowo-rekreacyjny" -> Unicode hyphen
$ws.Range("E199").Value = "ośrodek sportowo‐rekreacyjny"

# "ośrodek naukowo-badawczy" -> Unicode hyphen
$ws.Range("E206").Value = "ośrodek naukowo‐badawczy"

# "twierdza nad forteca" -> "twierdza lub forteca"
$ws.Range("E214").Value = "twierdza lub forteca"

# "gmina miejsko-wiejska" -> Unicode hyphen
$ws.Range("E231").Value = "gmina miejsko‐wiejska"

# "miasto w gminie miejsko-wiejskiej" -> Unicode hyphen
$ws.Range("E232").Value = "miasto w gminie miejsko‐wiejskiej"

# "obszar wiejski w gminie miejsko-wiejskiej" -> Unicode hyphen
$ws.Range("E233").Value = "obszar wiejski w gminie miejsko‐wiejskiej"

# --- Formatting touch-ups ---------------------------------------------------
# Re-typing/correcting the Polish labels above left those cells with an
# explicit black, 11pt Calibri font and a neutral (bottom/non-wrapping)
# alignment instead of the sheet's inherited default formatting. Re-apply the
# same explicit font to each corrected cell (one at a time, so the engine
# reuses a single shared style for all of them, matching how Excel collapses
# identical formatting into one cell style).
$editedCells = @("E50","E51","E72","E79","E122","E130","E158","E199","E206","E214","E231","E232","E233")
foreach ($addr in $editedCells) {
    $ws.Range($addr).Font.Name = "Calibri"
    $ws.Range($addr).Font.Size = 11
    $ws.Range($addr).Font.Color = 0
    $ws.Range($addr).VerticalAlignment = -4107
    $ws.Range($addr).WrapText = $false
    $ws.Range($addr).ShrinkToFit = $false
}

# F128's English text was edited too, and the cell's style was normalised to
# the same style already used by its row neighbour E128 - copy that cell's
# formatting across instead of inventing a new one.
$ws.Range("E128").Copy() | Out-Null
$ws.Range("F128").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
